# Updated cryptos list on Wed Sep 27 20:51:52 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the cryptos
# table, and swaps the Maker/Filecoin rows (32/33) to match the new
# ranking order.
#
# Note: some new Price values (e.g. "211.47") are plain decimals that
# Excel would otherwise auto-convert to a number on assignment. A
# leading apostrophe is used for those so the cell keeps storing text,
# exactly like typing '211.47 into a "General" formatted cell does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.264.38'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.594.99'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''211.47'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '''18.95'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").Value = '''0.0853'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '1.820.51'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.604.89'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("D16").Value = '''63.51'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '26.252.63'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '''229.85'
$ws.Range("E18").Value = '  +7.39%  '
$ws.Range("D19").Value = '0.0₃0720'
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = '''7.60'
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").Value = '''2.16'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("D24").Value = '''8.91'
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("D25").Value = '''145.99'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''6.98'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D29").Value = '''15.33'
$ws.Range("E29").Value = '  +1.52%  '
$ws.Range("D30").Value = '''0.0492'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.19'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.469.18'
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("D34").Value = '''2.93'
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''0.564'
$ws.Range("E37").Value = '  -4.14%  '
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '''5.73'
$ws.Range("E40").Value = '  -3.14%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").Value = '''2.16'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '''0.932'
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("D44").Value = '1.732.72'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").Value = '''0.752'
$ws.Range("E45").Value = '  -1.52%  '
$ws.Range("D46").Value = '''60.39'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").Value = '''87.76'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").Value = '''0.0502'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '''0.0946'
$ws.Range("E51").Value = '  -2.20%  '
